$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.563.56'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('D3').Value = '1.914.10'
$ws.Range('E3').Value = '  +5.48%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.76'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5181'
$ws.Range('E7').Value = '  +3.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3963'
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09705'
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.158'
$ws.Range('E10').Value = '  +4.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.94'
$ws.Range('E11').Value = '  +2.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.542'
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.24'
$ws.Range('E13').Value = '  +3.50%  '
$ws.Range('D14').Value = '1.917.53'
$ws.Range('E14').Value = '  +6.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.593'
$ws.Range('E15').Value = '  +4.40%  '
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001137'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.95'
$ws.Range('E18').Value = '  +1.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06659'
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.15'
$ws.Range('E20').Value = '  +5.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.321'
$ws.Range('E22').Value = '  +6.52%  '
$ws.Range('D23').Value = '28.626.59'
$ws.Range('E23').Value = '  +1.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.46'
$ws.Range('E24').Value = '  +3.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.299'
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.698'
$ws.Range('E26').Value = '  +12.07%  '
$ws.Range('B27').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C27').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D27').Value = '2.134.41'
$ws.Range('E27').Value = '  +5.66%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.24'
$ws.Range('E28').Value = '  +3.07%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '159.89'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.96'
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.105'
$ws.Range('E31').Value = '  +6.76%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1085'
$ws.Range('E32').Value = '  +1.99%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.750'
$ws.Range('E33').Value = '  +3.13%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.645'
$ws.Range('E34').Value = '  +1.33%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.911'
$ws.Range('E35').Value = '  +11.06%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06787'
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02435'
$ws.Range('E37').Value = '  +4.12%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.265'
$ws.Range('E38').Value = '  +7.53%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2227'
$ws.Range('E39').Value = '  +3.95%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.87'
$ws.Range('E40').Value = '  +5.08%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.117'
$ws.Range('E41').Value = '  +3.42%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6431'
$ws.Range('E42').Value = '  +3.85%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.189'
$ws.Range('E43').Value = '  +1.25%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.59'
$ws.Range('E45').Value = '  +3.15%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6085'
$ws.Range('E46').Value = '  +3.11%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.775'
$ws.Range('E47').Value = '  +2.22%  '
$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.281'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.034'
$ws.Range('E49').Value = '  +5.44%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.15'
$ws.Range('E50').Value = '  +0.73%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.213'
$ws.Range('E51').Value = '  +3.04%  '
